$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K12").Value = 13.76976495726495
$ws.Range("R12").Value = 1.867772014163364
$ws.Range("S12").Value = 2.045100507661769
$ws.Range("K13").Value = 13.76976495726495
$ws.Range("R13").Value = 3.202698560003336
$ws.Range("S13").Value = 3.918847414586112
$ws.Range("K26").Value = 21.19907407407406
$ws.Range("R26").Value = 3.745925377867032
$ws.Range("S26").Value = 4.780194493904943
$ws.Range("K27").Value = 21.19907407407406
$ws.Range("R27").Value = 2.015705049109126
$ws.Range("S27").Value = 2.22640738080769
$ws.Range("K28").Value = 2.356481481481501
$ws.Range("R28").Value = 1.678525338046114
$ws.Range("S28").Value = 1.817698795724144
$ws.Range("K29").Value = 2.356481481481501
$ws.Range("R29").Value = 2.619185573867416
$ws.Range("S29").Value = 3.069228739776626
$ws.Range("K36").Value = 20.68981481481483
$ws.Range("R36").Value = 3.702873262981926
$ws.Range("S36").Value = 4.709243017136692
$ws.Range("K37").Value = 20.68981481481483
$ws.Range("R37").Value = 2.004820578273036
$ws.Range("S37").Value = 2.212959200483225
$ws.Range("K38").Value = 21.28240740740739
$ws.Range("R38").Value = 2.017497406510892
$ws.Range("S38").Value = 2.228623569098047
$ws.Range("K39").Value = 21.28240740740739
$ws.Range("R39").Value = 3.753065762999045
$ws.Range("S39").Value = 4.792008787587529
$ws.Range("K40").Value = 15.74228395061728
$ws.Range("R40").Value = 3.33095021773865
$ws.Range("S40").Value = 4.115751405322535
$ws.Range("K41").Value = 15.74228395061728
$ws.Range("R41").Value = 1.904889690449167
$ws.Range("S41").Value = 2.090295475371289
$ws.Range("K42").Value = -1.819444444444444
$ws.Range("R42").Value = 2.455497817501559
$ws.Range("S42").Value = 2.843656807626497
$ws.Range("K43").Value = -1.819444444444444
$ws.Range("R43").Value = 1.618523362263702
$ws.Range("S43").Value = 1.746638928617865
$ws.Range("K46").Value = 13.46442495126706
$ws.Range("R46").Value = 1.862155209238257
$ws.Range("S46").Value = 2.038278558917324
$ws.Range("K47").Value = 13.46442495126706
$ws.Range("R47").Value = 3.183723085360894
$ws.Range("S47").Value = 3.890038835634921
$ws.Range("K54").Value = 19.79629629629628
$ws.Range("R54").Value = 3.629680458828347
$ws.Range("S54").Value = 4.589715938979482
$ws.Range("K55").Value = 19.79629629629628
$ws.Range("K56").Value = 19.79629629629628
$ws.Range("R56").Value = 1.98600466835246
$ws.Range("S56").Value = 2.18975222777657
$ws.Range("K59").Value = 21.79166666666666
$ws.Range("R59").Value = 2.028520339740724
$ws.Range("S59").Value = 2.242263395092639
$ws.Range("K60").Value = 21.79166666666666
$ws.Range("R60").Value = 3.797299903567984
$ws.Range("S60").Value = 4.865495608531995
$ws.Range("K61").Value = 12.51681286549706
$ws.Range("R61").Value = 1.844936767548521
$ws.Range("S61").Value = 2.017393709936214
$ws.Range("K62").Value = 12.51681286549706
$ws.Range("R62").Value = 3.126239257907711
$ws.Range("S62").Value = 3.803269267167952
$ws.Range("K63").Value = 15.74228395061728
$ws.Range("R63").Value = 1.904889690449167
$ws.Range("S63").Value = 2.090295475371289
$ws.Range("K64").Value = 15.74228395061728
$ws.Range("R64").Value = 3.33095021773865
$ws.Range("S64").Value = 4.115751405322535
$ws.Range("K65").Value = 15.74228395061728
$ws.Range("K71").Value = 13.75752314814816
$ws.Range("R71").Value = 1.867546171126113
$ws.Range("S71").Value = 2.044826120875009
$ws.Range("K72").Value = 13.75752314814816
$ws.Range("R72").Value = 3.201933436480062
$ws.Range("S72").Value = 3.917684201664166
$ws.Range("K73").Value = 14.96875
$ws.Range("K74").Value = 14.96875
$ws.Range("R74").Value = 1.890159325210871
$ws.Range("S74").Value = 2.072335994446373
$ws.Range("K75").Value = 14.96875
$ws.Range("R75").Value = 3.27945034353529
$ws.Range("S75").Value = 4.036221701795472
$ws.Range("K93").Value = 19.48611111111111
$ws.Range("R93").Value = 3.604943545926152
$ws.Range("S93").Value = 4.549628470864294
$ws.Range("K94").Value = 19.48611111111111
$ws.Range("R94").Value = 1.979555038534245
$ws.Range("S94").Value = 2.181809322722105
$ws.Range("K95").Value = 14.47727272727272
$ws.Range("R95").Value = 1.880917929007461
$ws.Range("S95").Value = 2.06108460959076
$ws.Range("K96").Value = 14.47727272727272
$ws.Range("R96").Value = 3.247548358074673
$ws.Range("S96").Value = 3.987268648345484
